$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 230, shifting existing rows 230..304 down to 231..305
$ws.Rows("230:230").Insert()

# Populate the newly inserted row 230 with the new data record
$ws.Range("A230").Value = 10
$ws.Range("B230").Value = "Vega Modelo de Temuco"
$ws.Range("C230").Value = "La Araucanía"
$ws.Range("D230").Value = 44924
$ws.Range("E230").Value = 9
$ws.Range("F230").Value = "Fruta"
$ws.Range("G230").Value = 100101
$ws.Range("H230").Value = "Berries"
$ws.Range("I230").Value = 100112025
$ws.Range("J230").Value = "Frutilla"
$ws.Range("K230").Value = "Sin especificar"
$ws.Range("L230").Value = "Primera"
$ws.Range("M230").Value = 135
$ws.Range("N230").Value = 6000
$ws.Range("O230").Value = 7000
$ws.Range("P230").Value = 6407
$ws.Range("Q230").Value = "$/caja 7 kilos"
$ws.Range("R230").Value = "Región de La Araucanía"
$ws.Range("S230").Value = 915
$ws.Range("T230").Value = 7
